# "Fruta / hortaliza, semanal" — weekly refresh of the Cebollín price series.
# A new weekly observation is prepended to the data block (row 89, right after
# the header + preceding records), pushing all existing records (old rows
# 89-136) down by one row to become rows 90-137.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 89; this shifts rows 89:136 down to 90:137 and
# inherits the number formatting (e.g. the date style on column D) from the
# row above, just like Excel's normal "Insert Sheet Rows" behavior.
$ws.Rows.Item(89).Insert()

# Populate the newly inserted row with the latest weekly record.
$ws.Cells.Item(89, 1).Value = 8
$ws.Cells.Item(89, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(89, 3).Value = "Coquimbo"
$ws.Cells.Item(89, 4).Value = 44529
$ws.Cells.Item(89, 5).Value = 4
$ws.Cells.Item(89, 6).Value = 100112037
$ws.Cells.Item(89, 7).Value = "Cebollín"
$ws.Cells.Item(89, 8).Value = "Sin especificar"
$ws.Cells.Item(89, 9).Value = "Primera"
$ws.Cells.Item(89, 10).Value = 3000
$ws.Cells.Item(89, 11).Value = 900
$ws.Cells.Item(89, 12).Value = 1000
$ws.Cells.Item(89, 13).Value = 950
$ws.Cells.Item(89, 14).Value = "`$/paquete 6 unidades"
$ws.Cells.Item(89, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(89, 16).Value = 158
$ws.Cells.Item(89, 17).Value = 6
$ws.Cells.Item(89, 18).Value = "Hortaliza"
